$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new columns M, N, O by copying formatting from column L ---
$ws.Range("L1:L7").Copy($ws.Range("M1:M7"))
$ws.Range("L1:L7").Copy($ws.Range("N1:N7"))
$ws.Range("L1:L7").Copy($ws.Range("O1:O7"))

# --- Row 1: headers ---
$ws.Range("M1").Value = "Q1 FY26"
$ws.Range("N1").Value = "Q2 FY26"
$ws.Range("O1").Value = "Q3 FY26"

# --- Row 2: report dates (update existing B:L, add M:O) ---
$ws.Range("B2").Value = 44773
$ws.Range("C2").Value = 44864
$ws.Range("D2").Value = 44955
$ws.Range("E2").Value = 45046
$ws.Range("F2").Value = 45137
$ws.Range("G2").Value = 45228
$ws.Range("H2").Value = 45319
$ws.Range("I2").Value = 45410
$ws.Range("J2").Value = 45501
$ws.Range("K2").Value = 45592
$ws.Range("L2").Value = 45683
$ws.Range("M2").Value = 45774
$ws.Range("N2").Value = 45865
$ws.Range("O2").Value = 45956

# --- Row 3: Data Center (add M:O) ---
$ws.Range("M3").Value = 39112
$ws.Range("N3").Value = 41096
$ws.Range("O3").Value = 51215

# --- Row 4: Gaming (add M:O) ---
$ws.Range("M4").Value = 3763
$ws.Range("N4").Value = 4287
$ws.Range("O4").Value = 4265

# --- Row 5: Professional Visualization (add M:O) ---
$ws.Range("M5").Value = 509
$ws.Range("N5").Value = 601
$ws.Range("O5").Value = 760

# --- Row 6: Automotive (add M:O) ---
$ws.Range("M6").Value = 567
$ws.Range("N6").Value = 586
$ws.Range("O6").Value = 592

# --- Row 7: OEM & Other (add M:O) ---
$ws.Range("M7").Value = 111
$ws.Range("N7").Value = 173
$ws.Range("O7").Value = 174

# --- Column widths for new columns (approximate bestFit; engine quantizes
#     ColumnWidth to 1/6-pixel steps, so these inputs are chosen to land the
#     closest achievable output to the target 9.08984375 / 10.08984375) ---
$ws.Range("M1:N7").ColumnWidth = 8.1
$ws.Range("O1:O7").ColumnWidth = 9.1

# --- Selection state ---
[void]$ws.Range("O2").Select()
